$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SQL")

# Insert two new rows before the current last row (row 24), pushing the
# existing row 24 ("select * from OLE.PORTAL_USER p join OLE.portal_user_tin...")
# down to row 26. The newly inserted rows 24/25 inherit the A/B column
# formatting (styles 5/7) from row 23 automatically.
$ws.Rows("24:25").Insert()

# Row 24: new Sno "23" + new query text
$ws.Range("A24").Value = "'23"
$ws.Range("A23").Copy()
$ws.Range("A24").PasteSpecial(-4122)

$query = " select t.PROV_TIN_NBR  from OLE.PORTAL_USER_TIN t join `n OLE.PORTAL_USER u on t.PORTAL_USER_ID=u.PORTAL_USER_ID `n where u.SSO_ID='{`$id}' and  u.STS_CD='A'"
$ws.Range("B24").Value = $query
$ws.Range("B23").Copy()
$ws.Range("B24").PasteSpecial(-4122)

$ws.Rows(24).RowHeight = 43.2

# Row 25 stays blank (spacer row) - formatting already inherited from insert.

# Update the view: selection now sits on B24.
$ws.Range("B24").Select() | Out-Null
